$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E2"  = 16.38969999999999
    "C3"  = -11.383
    "C14" = -13.65519999999999
    "C21" = -13.2738
    "C23" = -12.58930000000001
    "C25" = -11.3514
    "E25" = 17.12690000000001
    "C26" = -12.64920000000001
    "E27" = 16.69299999999998
    "C29" = -10.85620000000001
    "E31" = 16.7698
    "E39" = 16.199
    "E48" = 17.459
    "E51" = 17.0467
    "E52" = 16.88120000000001
    "C53" = -10.35640000000001
    "E55" = 16.58140000000001
    "E56" = 16.3745
    "C57" = -14.11039999999999
    "E57" = 16.7891
    "C59" = -12.8144
    "C69" = -11.5109
    "E73" = 17.2075
    "C79" = -10.95740000000001
    "C83" = -13.7639
    "E89" = 17.45280000000001
    "E90" = 16.40659999999999
    "C91" = -10.3564
    "E92" = 18.59120000000002
    "C93" = -11.3952
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
